$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") rows 2-67 all move from serial date 45171 (2023-09-02)
# to serial date 45172 (2023-09-03).
$ws.Range("C2:C67").Value = 45172
